# "hide filter in locations lists"
#
# The sample's three header cells used the raw field names coming back
# from the API ("...Id"); rename them to the friendlier labels shown in
# the locations filter lists, and scroll/re-select the sheet so the
# leading (now-hidden) filter column isn't the first thing in view.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Shared-string header renames: TerminalId -> Terminal, ErrandTypeId ->
# ErrandType, AssigneeId -> Assignee.
$ws.Range("A1").Value = "Terminal"
$ws.Range("B1").Value = "ErrandType"
$ws.Range("C1").Value = "Assignee"

# Scroll the view one column to the right (so column A is hidden off the
# left edge) and move the active selection from C5 to C2.
$excel.Goto($ws.Range("C2"), $true)
